$d = $word.ActiveDocument

$d.Content.Find.Execute("686÷2=343, 0", $true, $false, $false, $false, $false, $true, 1, $false, "767÷5=153, 2", 2) | Out-Null
$d.Content.Find.Execute("278÷7=39, 5", $true, $false, $false, $false, $false, $true, 1, $false, "998÷2=499, 0", 2) | Out-Null
$d.Content.Find.Execute("660÷4=165, 0", $true, $false, $false, $false, $false, $true, 1, $false, "594÷9=66, 0", 2) | Out-Null
$d.Content.Find.Execute("279÷4=69, 3", $true, $false, $false, $false, $false, $true, 1, $false, "347÷3=115, 2", 2) | Out-Null
$d.Content.Find.Execute("399÷9=44, 3", $true, $false, $false, $false, $false, $true, 1, $false, "569÷9=63, 2", 2) | Out-Null
$d.Content.Find.Execute("401÷9=44, 5", $true, $false, $false, $false, $false, $true, 1, $false, "766÷4=191, 2", 2) | Out-Null
$d.Content.Find.Execute("496÷7=70, 6", $true, $false, $false, $false, $false, $true, 1, $false, "895÷8=111, 7", 2) | Out-Null
$d.Content.Find.Execute("694÷7=99, 1", $true, $false, $false, $false, $false, $true, 1, $false, "888÷3=296, 0", 2) | Out-Null
$d.Content.Find.Execute("430÷7=61, 3", $true, $false, $false, $false, $false, $true, 1, $false, "623÷3=207, 2", 2) | Out-Null
$d.Content.Find.Execute("524÷8=65, 4", $true, $false, $false, $false, $false, $true, 1, $false, "776÷5=155, 1", 2) | Out-Null
$d.Content.Find.Execute("169÷2=84, 1", $true, $false, $false, $false, $false, $true, 1, $false, "774÷3=258, 0", 2) | Out-Null
$d.Content.Find.Execute("737÷9=81, 8", $true, $false, $false, $false, $false, $true, 1, $false, "841÷6=140, 1", 2) | Out-Null
$d.Content.Find.Execute("186÷4=46, 2", $true, $false, $false, $false, $false, $true, 1, $false, "163÷6=27, 1", 2) | Out-Null
$d.Content.Find.Execute("430÷2=215, 0", $true, $false, $false, $false, $false, $true, 1, $false, "370÷9=41, 1", 2) | Out-Null
$d.Content.Find.Execute("314÷9=34, 8", $true, $false, $false, $false, $false, $true, 1, $false, "945÷4=236, 1", 2) | Out-Null
$d.Content.Find.Execute("847÷7=121, 0", $true, $false, $false, $false, $false, $true, 1, $false, "367÷2=183, 1", 2) | Out-Null
$d.Content.Find.Execute("808÷4=202, 0", $true, $false, $false, $false, $false, $true, 1, $false, "838÷5=167, 3", 2) | Out-Null
$d.Content.Find.Execute("318÷6=53, 0", $true, $false, $false, $false, $false, $true, 1, $false, "773÷3=257, 2", 2) | Out-Null
$d.Content.Find.Execute("639÷2=319, 1", $true, $false, $false, $false, $false, $true, 1, $false, "110÷9=12, 2", 2) | Out-Null
$d.Content.Find.Execute("976÷7=139, 3", $true, $false, $false, $false, $false, $true, 1, $false, "197÷8=24, 5", 2) | Out-Null
$d.Content.Find.Execute("568÷5=113, 3", $true, $false, $false, $false, $false, $true, 1, $false, "429÷6=71, 3", 2) | Out-Null
$d.Content.Find.Execute("622÷2=311, 0", $true, $false, $false, $false, $false, $true, 1, $false, "337÷4=84, 1", 2) | Out-Null
$d.Content.Find.Execute("833÷8=104, 1", $true, $false, $false, $false, $false, $true, 1, $false, "394÷7=56, 2", 2) | Out-Null
$d.Content.Find.Execute("630÷4=157, 2", $true, $false, $false, $false, $false, $true, 1, $false, "208÷8=26, 0", 2) | Out-Null
$d.Content.Find.Execute("692÷5=138, 2", $true, $false, $false, $false, $false, $true, 1, $false, "126÷5=25, 1", 2) | Out-Null
